$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4500.6
$ws.Range("I98").Value = 3888.889
$ws.Range("K98").Value = 3888.889
$ws.Range("M98").Value = -2390.889

$ws.Range("H116").Value = 9617435
$ws.Range("I116").Value = 12822629
$ws.Range("J116").Value = 1853
$ws.Range("K116").Value = 12822629
$ws.Range("L116").Value = 1853
$ws.Range("M116").Value = -12819187
$ws.Range("N116").Value = -8737

$ws.Range("H122").Value = 4500.6
$ws.Range("I122").Value = 3888.889
$ws.Range("K122").Value = 11666.667
$ws.Range("M122").Value = -9216.667000000001

$ws.Range("H132").Value = 4103.3
$ws.Range("I132").Value = 3324.76
$ws.Range("J132").Value = 7996
$ws.Range("K132").Value = 9974.280000000001
$ws.Range("L132").Value = 23988
$ws.Range("M132").Value = -7444.280000000001
$ws.Range("N132").Value = -29048

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1389.4286
$ws.Range("I2").Value = 1400.6072
$ws.Range("J2").Value = 1344.7142
$ws.Range("K2").Value = 1400.6072
$ws.Range("L2").Value = 1344.7142
$ws.Range("M2").Value = -1287.6072
$ws.Range("N2").Value = -1570.7142

$ws.Range("H23").Value = 15999
$ws.Range("I23").Value = 14998.75
$ws.Range("J23").Value = 20000
$ws.Range("K23").Value = 14998.75
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = -14739.75
$ws.Range("N23").Value = -20518

$ws.Range("H27").Value = 7950
$ws.Range("J27").Value = 7950
$ws.Range("L27").Value = 7950
$ws.Range("N27").Value = -8318

$ws.Range("H32").Value = 6774.21
$ws.Range("I32").Value = 3061.169
$ws.Range("J32").Value = 19204.826
$ws.Range("K32").Value = 3061.169
$ws.Range("L32").Value = 19204.826
$ws.Range("M32").Value = -2774.169
$ws.Range("N32").Value = -19778.826

$ws.Range("H36").Value = 9000
$ws.Range("J36").Value = 10250
$ws.Range("L36").Value = 10250
$ws.Range("N36").Value = -10942

$ws.Range("H61").Value = 2502.6667
$ws.Range("I61").Value = 2400.4
$ws.Range("J61").Value = 3014
$ws.Range("K61").Value = 2400.4
$ws.Range("L61").Value = 3014
$ws.Range("M61").Value = -2188.4
$ws.Range("N61").Value = -3438

$ws.Range("H116").Value = 1389.4286
$ws.Range("I116").Value = 1400.6072
$ws.Range("J116").Value = 1344.7142
$ws.Range("K116").Value = 1400.6072
$ws.Range("L116").Value = 1344.7142
$ws.Range("M116").Value = 893.3928000000001
$ws.Range("N116").Value = -5932.7142

$ws.Range("H136").Value = 2502.6667
$ws.Range("I136").Value = 2400.4
$ws.Range("J136").Value = 3014
$ws.Range("K136").Value = 7201.200000000001
$ws.Range("L136").Value = 9042
$ws.Range("M136").Value = -4651.200000000001
$ws.Range("N136").Value = -14142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1389.4286
$ws.Range("I3").Value = 1400.6072
$ws.Range("J3").Value = 1344.7142
$ws.Range("K3").Value = 1400.6072
$ws.Range("L3").Value = 1344.7142
$ws.Range("M3").Value = -1286.6072
$ws.Range("N3").Value = -1572.7142

$ws.Range("H33").Value = 16382.857
$ws.Range("I33").Value = 8400
$ws.Range("J33").Value = 19576
$ws.Range("K33").Value = 8400
$ws.Range("L33").Value = 19576
$ws.Range("M33").Value = -8064
$ws.Range("N33").Value = -20248

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H134").Value = 2438.6
$ws.Range("I134").Value = 2138.8572
$ws.Range("J134").Value = 3138
$ws.Range("K134").Value = 6416.571599999999
$ws.Range("L134").Value = 9414
$ws.Range("M134").Value = -3881.571599999999
$ws.Range("N134").Value = -14484

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2593.7612
$ws.Range("I31").Value = 1717.5385
$ws.Range("K31").Value = 1717.5385
$ws.Range("M31").Value = -1422.5385

$ws.Range("H34").Value = 2593.7612
$ws.Range("I34").Value = 1717.5385
$ws.Range("K34").Value = 1717.5385
$ws.Range("M34").Value = -1515.5385

$ws.Range("H58").Value = 3217.1667
$ws.Range("I58").Value = 1636
$ws.Range("K58").Value = 1636
$ws.Range("M58").Value = -1433

$ws.Range("H106").Value = 34548.855
$ws.Range("J106").Value = 34548.855
$ws.Range("L106").Value = 34548.855
$ws.Range("N106").Value = -37072.855

$ws.Range("H136").Value = 3217.1667
$ws.Range("I136").Value = 1636
$ws.Range("K136").Value = 4908
$ws.Range("M136").Value = -2358

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5115
$ws.Range("I3").Value = 3402.5
$ws.Range("J3").Value = 5971.25
$ws.Range("K3").Value = 10207.5
$ws.Range("L3").Value = 17913.75
$ws.Range("M3").Value = -10095.5
$ws.Range("N3").Value = -18137.75

$ws.Range("H114").Value = 7188.75
$ws.Range("I114").Value = 6585.3335
$ws.Range("J114").Value = 8999
$ws.Range("K114").Value = 19756.0005
$ws.Range("L114").Value = 26997
$ws.Range("M114").Value = -16502.0005
$ws.Range("N114").Value = -33505

$ws.Range("H140").Value = 2714043
$ws.Range("I140").Value = 3584346
$ws.Range("K140").Value = 10753038
$ws.Range("M140").Value = -10747858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7241.4
$ws.Range("I70").Value = 8377.5
$ws.Range("J70").Value = 4590.5
$ws.Range("K70").Value = 8377.5
$ws.Range("L70").Value = 4590.5
$ws.Range("M70").Value = -8107.5
$ws.Range("N70").Value = -5130.5

$ws.Range("H73").Value = 7241.4
$ws.Range("I73").Value = 8377.5
$ws.Range("J73").Value = 4590.5
$ws.Range("K73").Value = 8377.5
$ws.Range("L73").Value = 4590.5
$ws.Range("M73").Value = -7441.5
$ws.Range("N73").Value = -6462.5

$ws.Range("H98").Value = 23000
$ws.Range("J98").Value = 23000
$ws.Range("L98").Value = 23000
$ws.Range("N98").Value = -28990

$ws.Range("H132").Value = 5839.1177
$ws.Range("I132").Value = 7049
$ws.Range("J132").Value = 4478
$ws.Range("K132").Value = 21147
$ws.Range("L132").Value = 13434
$ws.Range("M132").Value = -18617
$ws.Range("N132").Value = -18494

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2860.5134
$ws.Range("I136").Value = 1670.3914
$ws.Range("J136").Value = 4815.7144
$ws.Range("K136").Value = 5011.174199999999
$ws.Range("L136").Value = 14447.1432
$ws.Range("M136").Value = -2461.174199999999
$ws.Range("N136").Value = -19547.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 43166.668
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 43166.668
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 43166.668
$ws.Range("N105").Value = -50154.668
$ws.Range("M105").ClearContents()

$ws.Range("H136").Value = 3499.6047
$ws.Range("I136").Value = 3403
$ws.Range("J136").Value = 3922.25
$ws.Range("K136").Value = 10209
$ws.Range("L136").Value = 11766.75
$ws.Range("M136").Value = -7659
$ws.Range("N136").Value = -16866.75
